$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1250
$ws.Range("I20").Value = 1250
$ws.Range("K20").Value = 1250
$ws.Range("M20").Value = -1020

$ws.Range("H35").Value = 1250
$ws.Range("I35").Value = 1250
$ws.Range("K35").Value = 1250
$ws.Range("M35").Value = -871

$ws.Range("H42").Value = 2273015.8
$ws.Range("I42").Value = 3125095.5
$ws.Range("J42").Value = 802.6667
$ws.Range("K42").Value = 9375286.5
$ws.Range("L42").Value = 2408.0001
$ws.Range("M42").Value = -9375056.5
$ws.Range("N42").Value = -2868.0001

$ws.Range("H62").Value = 2588.9092
$ws.Range("I62").Value = 2050
$ws.Range("J62").Value = 2708.6667
$ws.Range("K62").Value = 2050
$ws.Range("L62").Value = 2708.6667
$ws.Range("M62").Value = -1426
$ws.Range("N62").Value = -3956.6667

$ws.Range("H65").Value = 2588.9092
$ws.Range("I65").Value = 2050
$ws.Range("J65").Value = 2708.6667
$ws.Range("K65").Value = 10250
$ws.Range("L65").Value = 13543.3335
$ws.Range("M65").Value = -7130
$ws.Range("N65").Value = -19783.3335

$ws.Range("H96").Value = 1249.75
$ws.Range("J96").Value = 1833
$ws.Range("L96").Value = 5499
$ws.Range("N96").Value = -8245

$ws.Range("H103").Value = 383
$ws.Range("I103").Value = 475
$ws.Range("J103").Value = 199
$ws.Range("K103").Value = 1425
$ws.Range("L103").Value = 597
$ws.Range("M103").Value = -839
$ws.Range("N103").Value = -1769

$ws.Range("H111").Value = 9633.154
$ws.Range("I111").Value = 11781.9
$ws.Range("J111").Value = 2470.6667
$ws.Range("K111").Value = 35345.7
$ws.Range("L111").Value = 7412.000100000001
$ws.Range("M111").Value = -32278.7
$ws.Range("N111").Value = -13546.0001

$ws.Range("H116").Value = 2500
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H129").Value = 1998.8085
$ws.Range("I129").Value = 8169.5386
$ws.Range("J129").Value = 1008.44446
$ws.Range("K129").Value = 24508.6158
$ws.Range("L129").Value = 3025.33338
$ws.Range("M129").Value = -19508.6158
$ws.Range("N129").Value = -13025.33338

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 38090.965
$ws.Range("I2").Value = 1216.3889
$ws.Range("K2").Value = 1216.3889
$ws.Range("M2").Value = -1103.3889

$ws.Range("H3").Value = 8000
$ws.Range("J3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("N3").Value = -8230

$ws.Range("H6").Value = 8690.200000000001
$ws.Range("I6").Value = 13025.5
$ws.Range("J6").Value = 5800
$ws.Range("K6").Value = 13025.5
$ws.Range("L6").Value = 5800
$ws.Range("M6").Value = -12852.5
$ws.Range("N6").Value = -6146

$ws.Range("H69").Value = 28125.916
$ws.Range("J69").Value = 28125.916
$ws.Range("L69").Value = 28125.916
$ws.Range("N69").Value = -29623.916

$ws.Range("H72").Value = 28125.916
$ws.Range("J72").Value = 28125.916
$ws.Range("L72").Value = 84377.74800000001
$ws.Range("N72").Value = -91865.74800000001

$ws.Range("H102").Value = 64450.812
$ws.Range("I102").Value = 252370
$ws.Range("K102").Value = 252370
$ws.Range("M102").Value = -250748

$ws.Range("H116").Value = 38090.965
$ws.Range("I116").Value = 1216.3889
$ws.Range("K116").Value = 1216.3889
$ws.Range("M116").Value = 1077.6111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 38090.965
$ws.Range("I3").Value = 1216.3889
$ws.Range("K3").Value = 1216.3889
$ws.Range("M3").Value = -1102.3889

$ws.Range("H8").Value = 43900
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 43900
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 43900
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -44180

$ws.Range("H70").Value = 27292.584
$ws.Range("J70").Value = 27292.584
$ws.Range("L70").Value = 27292.584
$ws.Range("N70").Value = -27878.584

$ws.Range("H73").Value = 27292.584
$ws.Range("J73").Value = 27292.584
$ws.Range("L73").Value = 27292.584
$ws.Range("N73").Value = -29320.584

$ws.Range("H99").Value = 1918
$ws.Range("I99").Value = 3580
$ws.Range("J99").Value = 1830.5264
$ws.Range("K99").Value = 3580
$ws.Range("L99").Value = 1830.5264
$ws.Range("M99").Value = -2082
$ws.Range("N99").Value = -4826.5264

$ws.Range("H134").Value = 2636.1853
$ws.Range("I134").Value = 2368.15
$ws.Range("J134").Value = 3402
$ws.Range("K134").Value = 7104.450000000001
$ws.Range("L134").Value = 10206
$ws.Range("M134").Value = -4569.450000000001
$ws.Range("N134").Value = -15276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H134").Value = 1970.3334
$ws.Range("I134").Value = 1864.4
$ws.Range("K134").Value = 5593.200000000001
$ws.Range("M134").Value = -3058.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 126.2
$ws.Range("I6").Value = 118
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 354
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -241
$ws.Range("N6").Value = -826

$ws.Range("H21").Value = 91
$ws.Range("I21").Value = 91
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 273
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -100
$ws.Range("N21").ClearContents()

$ws.Range("H131").Value = 786.17
$ws.Range("J131").Value = 802.2371000000001
$ws.Range("L131").Value = 2406.7113
$ws.Range("N131").Value = -12486.7113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5000000
$ws.Range("I7").Value = 5000000
$ws.Range("K7").Value = 5000000
$ws.Range("M7").Value = -4999888

$ws.Range("H8").Value = 5000000
$ws.Range("I8").Value = 5000000
$ws.Range("K8").Value = 5000000
$ws.Range("M8").Value = -4999861

$ws.Range("H11").Value = 5292427
$ws.Range("I11").Value = 5622613.5
$ws.Range("K11").Value = 5622613.5
$ws.Range("M11").Value = -5622474.5

$ws.Range("H99").Value = 3261.5715
$ws.Range("I99").Value = 3261.5715
$ws.Range("K99").Value = 3261.5715
$ws.Range("M99").Value = -1015.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1543.1765
$ws.Range("J61").Value = 1388.5714
$ws.Range("L61").Value = 1388.5714
$ws.Range("N61").Value = -1792.5714

$ws.Range("H100").Value = 2497.6667
$ws.Range("I100").Value = 2171.2856
$ws.Range("J100").Value = 3640
$ws.Range("K100").Value = 2171.2856
$ws.Range("L100").Value = 3640
$ws.Range("M100").Value = -1630.2856
$ws.Range("N100").Value = -4722

$ws.Range("H113").Value = 1543.1765
$ws.Range("J113").Value = 1388.5714
$ws.Range("L113").Value = 1388.5714
$ws.Range("N113").Value = -5728.5714

$ws.Range("H132").Value = 4425.3184
$ws.Range("I132").Value = 4292.9
$ws.Range("K132").Value = 12878.7
$ws.Range("M132").Value = -10348.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 16335
$ws.Range("I13").Value = 502.5
$ws.Range("K13").Value = 502.5
$ws.Range("M13").Value = -362.5

$ws.Range("H100").Value = 501750
$ws.Range("I100").Value = 1000000
$ws.Range("K100").Value = 2000000
$ws.Range("M100").Value = -1999459

$ws.Range("H107").Value = 200772
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 200772
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 602316
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -606156

$ws.Range("H113").Value = 616.0968
$ws.Range("I113").Value = 431.85
$ws.Range("K113").Value = 1295.55
$ws.Range("M113").Value = 874.4499999999998

$ws.Range("H132").Value = 2139
$ws.Range("I132").Value = 2033.6666
$ws.Range("J132").Value = 2549.8
$ws.Range("K132").Value = 6100.9998
$ws.Range("L132").Value = 7649.400000000001
$ws.Range("M132").Value = -3570.9998
$ws.Range("N132").Value = -12709.4

Write-Output "Applied all cell updates."